# "Generate Report for Handback" -- append the handback row for the
# 9200d23d-3e7e-43e6-82f7-3a8b98bae92c file to the Overview sheet and to
# each per-locale (zh-cn / de-de) detail sheet, growing each sheet's table
# by one row (row 3) and wiring up the matching hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ovTable = $ov.ListObjects.Item(1)
$ovTable.ListRows.Add() | Out-Null

$ov.Range("A3").Value = "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md"

$ov.Range("B3").Value = "e2e\9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md"
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/839bae1745b299482fac4f1190bdee217237ff58/e2e/9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md") | Out-Null
$ov.Range("B3").Style = "Hyperlink"

$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"

$ov.Range("G3").Value = "2016-10-25 02:15:23"
$ov.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zhTable = $zh.ListObjects.Item(1)
$zhTable.ListRows.Add() | Out-Null

$zh.Range("A3").Value = "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md"
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c7396d51095e0cee356cf74f1ad5cddf5f5228b6/e2e/9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md", [System.Type]::Missing, [System.Type]::Missing, "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md") | Out-Null
$zh.Range("A3").Style = "Hyperlink"

$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "ht"
$zh.Range("F3").Value = "True"

$zh.Range("G3").Value = "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.5e86f675d6b1acf32b98a1bd7789d6d4a0f91a2c.zh-cn.xlf"

$zh.Range("H3").Value = "2016-10-25 02:15:10"
$zh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$zh.Range("I3").Value = "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md"
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c7396d51095e0cee356cf74f1ad5cddf5f5228b6/e2e/9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md", [System.Type]::Missing, [System.Type]::Missing, "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md") | Out-Null
$zh.Range("I3").Style = "Hyperlink"

$zh.Range("J3").Value = "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.5e86f675d6b1acf32b98a1bd7789d6d4a0f91a2c.zh-cn.xlf"

$zh.Range("K3").Value = "2016-10-25 02:15:52"
$zh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$zh.Range("L3").Value = ""
$zh.Range("M3").Value = "True"
$zh.Range("N3").Value = ""
$zh.Range("O3").Value = "False"
$zh.Range("P3").Value = ""

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$deTable = $de.ListObjects.Item(1)
$deTable.ListRows.Add() | Out-Null

$de.Range("A3").Value = "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md"
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7acadb984a87f03f312ba53d790ff7ee93bdfdd0/e2e/9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md", [System.Type]::Missing, [System.Type]::Missing, "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md") | Out-Null
$de.Range("A3").Style = "Hyperlink"

$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "ht"
$de.Range("F3").Value = "True"

$de.Range("G3").Value = "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.5e86f675d6b1acf32b98a1bd7789d6d4a0f91a2c.de-de.xlf"

$de.Range("H3").Value = "2016-10-25 02:16:11"
$de.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$de.Range("I3").Value = "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md"
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7acadb984a87f03f312ba53d790ff7ee93bdfdd0/e2e/9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md", [System.Type]::Missing, [System.Type]::Missing, "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.md") | Out-Null
$de.Range("I3").Style = "Hyperlink"

$de.Range("J3").Value = "9200d23d-3e7e-43e6-82f7-3a8b98bae92c.5e86f675d6b1acf32b98a1bd7789d6d4a0f91a2c.de-de.xlf"

$de.Range("K3").Value = "2016-10-25 02:16:11"
$de.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$de.Range("L3").Value = ""
$de.Range("M3").Value = "True"
$de.Range("N3").Value = ""
$de.Range("O3").Value = "False"
$de.Range("P3").Value = ""

Write-Output "Generated handback report rows for Overview, zh-cn, de-de"
